# etl_test.xlsx - "get_last_val now stops at blank"
#
# The Reference sheet's trailing helper column (Q) was a blank spacer
# between the lookup columns (…P) and the get_last_val/find_last_val_col
# test column (old R:Z). Removing that spacer shifts the test column from
# Z to Y, and a new sample value ("more") is added in the row that used to
# be a blank gap (old Z5 / new Y5) so the fixture actually exercises the
# "stop scanning at the first blank cell" behaviour.

$wb = $excel.ActiveWorkbook

$wsDb  = $wb.Worksheets.Item("Database")
$wsRef = $wb.Worksheets.Item("Reference")

# --- Reference sheet: drop the blank spacer column Q (shifts R:Z -> Q:Y) ---
$wsRef.Columns("Q:Q").Delete() | Out-Null

# --- New fixture value in what is now the blank row of the "last value" column ---
$wsRef.Range("Y5").Value = "more"

# --- Header rows grow a couple of points taller on both sheets ---
$wsDb.Rows.Item(1).RowHeight  = 66
$wsRef.Rows.Item(1).RowHeight = 57

# --- Selection/active-sheet state: Reference (with its new last column) is now
#     the active tab, selected at its new last cell; Database's old selection
#     moves down to A14 and is no longer the active tab ---
$wsDb.Range("A14").Select() | Out-Null
$wsRef.Select() | Out-Null
$wsRef.Range("Y1").Select() | Out-Null

# --- Workbook recalculation switched to manual ---
$excel.Calculation = -4135
